# ============================================================
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4"
#    (before "总计"), with the same column layout/styling as
#    the other quarterly fund-holding sheets, populated with
#    the Q1-2022 fund holdings.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet and
#    renumber the existing index column.
# ============================================================

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet, positioned between
# "2021-Q4" and "总计".
# ------------------------------------------------------------
$prevQuarter = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $prevQuarter)
$newSheet.Name = "2022-Q1"

# Clone the column styling (header style + index-column style)
# from the "2021-Q4" sheet so the new sheet matches the look of
# the other quarterly sheets.
$template = $wb.Worksheets.Item("2021-Q4")
$template.Range("A1:H9").Copy()
$newSheet.Range("A1:H9").PasteSpecial(-4122)
$template.Range("A2:H2").Copy()
$newSheet.Range("A10:H11").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "'360006"
$newSheet.Cells.Item(2, 3).Value = "光大保德信新增长混合"
$newSheet.Cells.Item(2, 4).Value = "'21.71"
$newSheet.Cells.Item(2, 5).Value = "'88.07"
$newSheet.Cells.Item(2, 6).Value = "'4.06"
$newSheet.Cells.Item(2, 7).Value = "'0.8814"
$newSheet.Cells.Item(2, 8).Value = 6

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "'011104"
$newSheet.Cells.Item(3, 3).Value = "光大保德信智能汽车主题股票"
$newSheet.Cells.Item(3, 4).Value = "'10.77"
$newSheet.Cells.Item(3, 5).Value = "'90.06"
$newSheet.Cells.Item(3, 6).Value = "'5.11"
$newSheet.Cells.Item(3, 7).Value = "'0.5503"
$newSheet.Cells.Item(3, 8).Value = 5

$newSheet.Cells.Item(4, 1).Value = 2
$newSheet.Cells.Item(4, 2).Value = "'001740"
$newSheet.Cells.Item(4, 3).Value = "光大保德信中国制造2025灵活配置混合"
$newSheet.Cells.Item(4, 4).Value = "'11.43"
$newSheet.Cells.Item(4, 5).Value = "'86.23"
$newSheet.Cells.Item(4, 6).Value = "'3.21"
$newSheet.Cells.Item(4, 7).Value = "'0.3669"
$newSheet.Cells.Item(4, 8).Value = 9

$newSheet.Cells.Item(5, 1).Value = 3
$newSheet.Cells.Item(5, 2).Value = "'008347"
$newSheet.Cells.Item(5, 3).Value = "中信建投价值甄选混合A"
$newSheet.Cells.Item(5, 4).Value = "'4.45"
$newSheet.Cells.Item(5, 5).Value = "'72.24"
$newSheet.Cells.Item(5, 6).Value = "'2.73"
$newSheet.Cells.Item(5, 7).Value = "'0.1215"
$newSheet.Cells.Item(5, 8).Value = 5

$newSheet.Cells.Item(6, 1).Value = 4
$newSheet.Cells.Item(6, 2).Value = "'003822"
$newSheet.Cells.Item(6, 3).Value = "中信建投行业轮换混合A"
$newSheet.Cells.Item(6, 4).Value = "'3.07"
$newSheet.Cells.Item(6, 5).Value = "'72.09"
$newSheet.Cells.Item(6, 6).Value = "'2.73"
$newSheet.Cells.Item(6, 7).Value = "'0.0838"
$newSheet.Cells.Item(6, 8).Value = 6

$newSheet.Cells.Item(7, 1).Value = 5
$newSheet.Cells.Item(7, 2).Value = "'007468"
$newSheet.Cells.Item(7, 3).Value = "中信建投策略精选混合A"
$newSheet.Cells.Item(7, 4).Value = "'0.94"
$newSheet.Cells.Item(7, 5).Value = "'78.13"
$newSheet.Cells.Item(7, 6).Value = "'3.00"
$newSheet.Cells.Item(7, 7).Value = "'0.0282"
$newSheet.Cells.Item(7, 8).Value = 5

$newSheet.Cells.Item(8, 1).Value = 6
$newSheet.Cells.Item(8, 2).Value = "'002630"
$newSheet.Cells.Item(8, 3).Value = "江信瑞福灵活配置混合A"
$newSheet.Cells.Item(8, 4).Value = "'0.52"
$newSheet.Cells.Item(8, 5).Value = "'43.17"
$newSheet.Cells.Item(8, 6).Value = "'4.58"
$newSheet.Cells.Item(8, 7).Value = "'0.0238"
$newSheet.Cells.Item(8, 8).Value = 2

$newSheet.Cells.Item(9, 1).Value = 7
$newSheet.Cells.Item(9, 2).Value = "'002631"
$newSheet.Cells.Item(9, 3).Value = "江信瑞福灵活配置混合C"
$newSheet.Cells.Item(9, 4).Value = "'0.50"
$newSheet.Cells.Item(9, 5).Value = "'43.17"
$newSheet.Cells.Item(9, 6).Value = "'4.58"
$newSheet.Cells.Item(9, 7).Value = "'0.0229"
$newSheet.Cells.Item(9, 8).Value = 2

$newSheet.Cells.Item(10, 1).Value = 8
$newSheet.Cells.Item(10, 2).Value = "'003823"
$newSheet.Cells.Item(10, 3).Value = "中信建投行业轮换混合C"
$newSheet.Cells.Item(10, 4).Value = "'0.64"
$newSheet.Cells.Item(10, 5).Value = "'72.09"
$newSheet.Cells.Item(10, 6).Value = "'2.73"
$newSheet.Cells.Item(10, 7).Value = "'0.0175"
$newSheet.Cells.Item(10, 8).Value = 6

$newSheet.Cells.Item(11, 1).Value = 9
$newSheet.Cells.Item(11, 2).Value = "'007469"
$newSheet.Cells.Item(11, 3).Value = "中信建投策略精选混合C"
$newSheet.Cells.Item(11, 4).Value = "'0.40"
$newSheet.Cells.Item(11, 5).Value = "'78.13"
$newSheet.Cells.Item(11, 6).Value = "'3.00"
$newSheet.Cells.Item(11, 7).Value = "'0.0120"
$newSheet.Cells.Item(11, 8).Value = 5


# ------------------------------------------------------------
# Step 2: insert the new "2022-Q1" row into the "总计" summary
# sheet, right after the header, and renumber the index column.
# ------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Clone the formatting of the row immediately below (which still
# carries the original row-2 formatting) into the newly inserted
# blank row.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q1"
$summary.Cells.Item(2, 3).Value = 10
$summary.Cells.Item(2, 4).Value = 2.11

# Renumber the existing rows' index column (0..4 -> 1..5)
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
